$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new value, and whether to force-as-text
$updates = @(
    [PSCustomObject]@{ Row = 2; Col = 4; Value = "'29.122.49"; ForceText = $true },
    [PSCustomObject]@{ Row = 2; Col = 5; Value = "  -0.98%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 3; Col = 4; Value = "'1.834.72"; ForceText = $true },
    [PSCustomObject]@{ Row = 3; Col = 5; Value = "  -0.98%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 4; Col = 4; Value = "'0.9986"; ForceText = $true },
    [PSCustomObject]@{ Row = 4; Col = 5; Value = "  -0.16%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 5; Col = 4; Value = "'240.55"; ForceText = $true },
    [PSCustomObject]@{ Row = 5; Col = 5; Value = "  -1.91%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 6; Col = 4; Value = "'0.6644"; ForceText = $true },
    [PSCustomObject]@{ Row = 6; Col = 5; Value = "  -4.04%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 7; Col = 4; Value = "'0.9997"; ForceText = $true },
    [PSCustomObject]@{ Row = 7; Col = 5; Value = "  -0.11%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 8; Col = 4; Value = "'0.2945"; ForceText = $true },
    [PSCustomObject]@{ Row = 8; Col = 5; Value = "  -3.80%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 9; Col = 4; Value = "'0.07349"; ForceText = $true },
    [PSCustomObject]@{ Row = 10; Col = 4; Value = "'22.71"; ForceText = $true },
    [PSCustomObject]@{ Row = 10; Col = 5; Value = "  -3.28%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 11; Col = 4; Value = "'0.07683"; ForceText = $true },
    [PSCustomObject]@{ Row = 11; Col = 5; Value = "  -0.95%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 12; Col = 4; Value = "'5.013"; ForceText = $true },
    [PSCustomObject]@{ Row = 12; Col = 5; Value = "  -2.52%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 13; Col = 4; Value = "'1.792.34"; ForceText = $true },
    [PSCustomObject]@{ Row = 13; Col = 5; Value = "  -3.31%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 14; Col = 4; Value = "'0.6728"; ForceText = $true },
    [PSCustomObject]@{ Row = 14; Col = 5; Value = "  -2.95%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 15; Col = 4; Value = "'86.02"; ForceText = $true },
    [PSCustomObject]@{ Row = 15; Col = 5; Value = "  -5.41%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 16; Col = 4; Value = "'6.199"; ForceText = $true },
    [PSCustomObject]@{ Row = 16; Col = 5; Value = "  -1.66%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 17; Col = 4; Value = "'0.000008209"; ForceText = $true },
    [PSCustomObject]@{ Row = 17; Col = 5; Value = "  -0.80%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 18; Col = 4; Value = "'28.787.63"; ForceText = $true },
    [PSCustomObject]@{ Row = 18; Col = 5; Value = "  -2.18%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 19; Col = 4; Value = "'227.85"; ForceText = $true },
    [PSCustomObject]@{ Row = 19; Col = 5; Value = "  -3.56%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 20; Col = 4; Value = "'12.49"; ForceText = $true },
    [PSCustomObject]@{ Row = 20; Col = 5; Value = "  -1.70%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 21; Col = 4; Value = "'0.9985"; ForceText = $true },
    [PSCustomObject]@{ Row = 21; Col = 5; Value = "  -0.19%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 22; Col = 4; Value = "'7.239"; ForceText = $true },
    [PSCustomObject]@{ Row = 22; Col = 5; Value = "  -5.38%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 23; Col = 4; Value = "'0.9997"; ForceText = $true },
    [PSCustomObject]@{ Row = 23; Col = 5; Value = "  -0.14%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 24; Col = 4; Value = "'160.29"; ForceText = $true },
    [PSCustomObject]@{ Row = 24; Col = 5; Value = "  +0.20%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 25; Col = 4; Value = "'8.685"; ForceText = $true },
    [PSCustomObject]@{ Row = 25; Col = 5; Value = "  -2.65%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 26; Col = 4; Value = "'0.1396"; ForceText = $true },
    [PSCustomObject]@{ Row = 26; Col = 5; Value = "  -5.59%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 27; Col = 4; Value = "'18.01"; ForceText = $true },
    [PSCustomObject]@{ Row = 27; Col = 5; Value = "  -1.11%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 28; Col = 5; Value = "  -1.71%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 29; Col = 4; Value = "'4.202"; ForceText = $true },
    [PSCustomObject]@{ Row = 29; Col = 5; Value = "  -0.98%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 30; Col = 4; Value = "'4.077"; ForceText = $true },
    [PSCustomObject]@{ Row = 31; Col = 4; Value = "'1.191"; ForceText = $true },
    [PSCustomObject]@{ Row = 31; Col = 5; Value = "  -0.84%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 32; Col = 4; Value = "'0.05356"; ForceText = $true },
    [PSCustomObject]@{ Row = 32; Col = 5; Value = "  +2.38%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 33; Col = 4; Value = "'0.7487"; ForceText = $true },
    [PSCustomObject]@{ Row = 33; Col = 5; Value = "  -3.14%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 34; Col = 4; Value = "'1.850"; ForceText = $true },
    [PSCustomObject]@{ Row = 34; Col = 5; Value = "  -0.87%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 35; Col = 4; Value = "'1.133"; ForceText = $true },
    [PSCustomObject]@{ Row = 35; Col = 5; Value = "  -1.03%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 36; Col = 4; Value = "'2.680"; ForceText = $true },
    [PSCustomObject]@{ Row = 36; Col = 5; Value = "  -0.29%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 37; Col = 4; Value = "'1.324.30"; ForceText = $true },
    [PSCustomObject]@{ Row = 37; Col = 5; Value = "  -0.26%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 38; Col = 4; Value = "'0.01804"; ForceText = $true },
    [PSCustomObject]@{ Row = 38; Col = 5; Value = "  -3.20%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 39; Col = 4; Value = "'2.715"; ForceText = $true },
    [PSCustomObject]@{ Row = 39; Col = 5; Value = "  -0.24%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 40; Col = 4; Value = "'0.9220"; ForceText = $true },
    [PSCustomObject]@{ Row = 40; Col = 5; Value = "  -2.03%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 41; Col = 4; Value = "'5.953"; ForceText = $true },
    [PSCustomObject]@{ Row = 41; Col = 5; Value = "  +2.52%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 42; Col = 4; Value = "'0.9979"; ForceText = $true },
    [PSCustomObject]@{ Row = 42; Col = 5; Value = "  -0.26%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 43; Col = 4; Value = "'103.28"; ForceText = $true },
    [PSCustomObject]@{ Row = 43; Col = 5; Value = "  -2.41%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 44; Col = 2; Value = "XinFinNetwork"; ForceText = $false },
    [PSCustomObject]@{ Row = 44; Col = 3; Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"; ForceText = $false },
    [PSCustomObject]@{ Row = 44; Col = 4; Value = "'0.07913"; ForceText = $true },
    [PSCustomObject]@{ Row = 44; Col = 5; Value = "  +14.77%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 45; Col = 2; Value = "BabyDogeCoin"; ForceText = $false },
    [PSCustomObject]@{ Row = 45; Col = 3; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; ForceText = $false },
    [PSCustomObject]@{ Row = 45; Col = 4; Value = "'0.00000000124"; ForceText = $true },
    [PSCustomObject]@{ Row = 45; Col = 5; Value = "  +0.01%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 46; Col = 4; Value = "'0.5163"; ForceText = $true },
    [PSCustomObject]@{ Row = 46; Col = 5; Value = "  -1.21%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 47; Col = 2; Value = "RocketPoolETH"; ForceText = $false },
    [PSCustomObject]@{ Row = 47; Col = 3; Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; ForceText = $false },
    [PSCustomObject]@{ Row = 47; Col = 4; Value = "'1.931.12"; ForceText = $true },
    [PSCustomObject]@{ Row = 47; Col = 5; Value = "  -3.41%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 48; Col = 2; Value = "Aave"; ForceText = $false },
    [PSCustomObject]@{ Row = 48; Col = 3; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; ForceText = $false },
    [PSCustomObject]@{ Row = 48; Col = 4; Value = "'63.62"; ForceText = $true },
    [PSCustomObject]@{ Row = 48; Col = 5; Value = "  +1.15%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 49; Col = 4; Value = "'1.751"; ForceText = $true },
    [PSCustomObject]@{ Row = 49; Col = 5; Value = "  -1.74%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 50; Col = 4; Value = "'9.283"; ForceText = $true },
    [PSCustomObject]@{ Row = 50; Col = 5; Value = "  -4.28%  "; ForceText = $false },
    [PSCustomObject]@{ Row = 51; Col = 4; Value = "'0.05927"; ForceText = $true },
    [PSCustomObject]@{ Row = 51; Col = 5; Value = "  -0.43%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.Value = $u.Value
    if ($u.ForceText) {
        $cell.Style = "Normal"
    }
}
